# Apply the changes described by the commit diff:
# 1. Update the cached "datetimeFigureOut" field text from 15/01/2025 to
#    21/01/2025 on the slide master and every slide layout's date
#    placeholder.
# 2. On slide 1, merge the trailing " " + "2025" runs of the
#    "WEST edition 2025" textbox into a single " 2025" run (keeping the
#    box's autosize height unchanged).
# 3. On slide 1, nudge the "Image 49" picture's horizontal offset from
#    4881548 EMU to 4848296 EMU (Top stays put).

$p = $ppt.ActivePresentation

$oldDate = "15/01/2025"
$newDate = "21/01/2025"

function Update-DatePlaceholder($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.HasTextFrame) {
            if ($shp.TextFrame.HasText) {
                $t = $shp.TextFrame.TextRange.Text
                if ($t -eq $oldDate) {
                    $shp.TextFrame.TextRange.Text = $newDate
                }
            }
        }
    }
}

# 1a. Slide master date placeholder.
Update-DatePlaceholder $p.SlideMaster.Shapes

# 1b. Every slide layout's date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DatePlaceholder $layouts.Item($i).Shapes
}

# 2 & 3. Slide 1 specific tweaks.
$s1 = $p.Slides.Item(1)
for ($k = 1; $k -le $s1.Shapes.Count; $k++) {
    $shp = $s1.Shapes.Item($k)

    if ($shp.Name -eq "ZoneTexte 43") {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "WEST edition 2025") {
            $origHeight = $shp.Height
            # Characters(13,5) spans the trailing " " + "2025" runs;
            # rewriting it merges them into a single run " 2025".
            $chars = $tr.Characters(13, 5)
            $chars.Text = " 2025"
            # Re-assert the autofit height so the spAutoFit recalculation
            # triggered by the text edit doesn't leave a stray diff.
            $shp.Height = $origHeight
        }
    }
    elseif ($shp.Name -eq "Image 49") {
        if ([math]::Round($shp.Left * 12700) -eq 4881548) {
            # Round to the same 4-decimal-place point precision PowerPoint
            # itself stores, so the re-emitted EMU value lands exactly on
            # 4848296 instead of off-by-one from float division noise.
            $shp.Left = [math]::Round(4848296 / 12700, 4)
        }
    }
}
